# Update the "Assets" sheet: rename/rework a handful of file-path asset
# rows (ExcelFilePath / Elevresor_ExcelFilePath / Grundskola_ExcelFilePath
# are consolidated into FilePath_Masterfile_Elevresor /
# FilePath_Masterfile_Grundskola), and shift the remaining file-path rows
# down to make room for the restored "Conversionfile_Elevresor" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

# Row 5: Elevresor_ExcelFilePath / ExcelFilePath -> FilePath_Masterfile_Elevresor
$ws.Range("A5").Value = "FilePath_Masterfile_Elevresor"
$ws.Range("B5").Value = "FilePath_Masterfile_Elevresor"

# Row 6: Grundskola_ExcelFilePath / FilePath_Grundskola -> FilePath_Masterfile_Grundskola
$ws.Range("A6").Value = "FilePath_Masterfile_Grundskola"
$ws.Range("B6").Value = "FilePath_Masterfile_Grundskola"

# Row 7 (Skolportalen_Template) is unchanged.

# Row 8: FilePath_Elevresor -> FilePath_Grundskola
$ws.Range("A8").Value = "FilePath_Grundskola"
$ws.Range("B8").Value = "FilePath_Grundskola"

# Row 9: Conversionfile_Elevresor / ConversionExcelFilepath -> FilePath_Elevresor
$ws.Range("A9").Value = "FilePath_Elevresor"
$ws.Range("B9").Value = "FilePath_Elevresor"

# Row 10 (previously empty): restore Conversionfile_Elevresor / ConversionExcelFilepath
$ws.Range("A10").Value = "Conversionfile_Elevresor"
$ws.Range("B10").Value = "ConversionExcelFilepath"

# Keep the trailing formatted-row extent one row taller, matching the sheet's
# extended used range.
$ws.Rows.Item(1001).RowHeight = 14.25

# Leave the selection where the edit session ended.
$ws.Range("B8").Select()
